{"js": "const replacements = [\n  [\"2025-11-06 Thursday\", \"2025-11-07 Friday\"],\n  [\"784\u00d72=\", \"575\u00d75=\"],\n  [\"289\u00d76=\", \"502\u00d75=\"],\n  [\"653\u00d78=\", \"380\u00d74=\"],\n  [\"159\u00d78=\", \"854\u00d77=\"],\n  [\"353\u00d73=\", \"416\u00d73=\"],\n  [\"798\u00d73=\", \"855\u00d72=\"],\n  [\"393\u00d77=\", \"731\u00d72=\"],\n  [\"552\u00d76=\", \"597\u00d72=\"],\n  [\"297\u00d72=\", \"313\u00d79=\"],\n  [\"379\u00d73=\", \"653\u00d75=\"],\n  [\"903\u00d74=\", \"512\u00d77=\"],\n  [\"880\u00d76=\", \"512\u00d79=\"],\n  [\"729\u00d79=\", \"822\u00d72=\"],\n  [\"825\u00d78=\", \"782\u00d77=\"],\n  [\"441\u00d78=\", \"815\u00d76=\"],\n  [\"858\u00d72=\", \"143\u00d72=\"],\n  [\"404\u00d74=\", \"658\u00d78=\"],\n  [\"737\u00d75=\", \"737\u00d78=\"],\n  [\"531\u00d76=\", \"618\u00d76=\"],\n  [\"660\u00d77=\", \"355\u00d72=\"],\n  [\"636\u00d75=\", \"994\u00d79=\"],\n  [\"910\u00d77=\", \"901\u00d73=\"],\n  [\"998\u00d78=\", \"650\u00d72=\"],\n  [\"198\u00d76=\", \"116\u00d78=\"],\n  [\"912\u00d74=\", \"583\u00d72=\"],\n];\n\nconst body = context.document.body;\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $null = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Text \"2025-11-06 Thursday\" \"2025-11-07 Friday\"\nReplace-Text \"784\u00d72=\" \"575\u00d75=\"\nReplace-Text \"289\u00d76=\" \"502\u00d75=\"\nReplace-Text \"653\u00d78=\" \"380\u00d74=\"\nReplace-Text \"159\u00d78=\" \"854\u00d77=\"\nReplace-Text \"353\u00d73=\" \"416\u00d73=\"\nReplace-Text \"798\u00d73=\" \"855\u00d72=\"\nReplace-Text \"393\u00d77=\" \"731\u00d72=\"\nReplace-Text \"552\u00d76=\" \"597\u00d72=\"\nReplace-Text \"297\u00d72=\" \"313\u00d79=\"\nReplace-Text \"379\u00d73=\" \"653\u00d75=\"\nReplace-Text \"903\u00d74=\" \"512\u00d77=\"\nReplace-Text \"880\u00d76=\" \"512\u00d79=\"\nReplace-Text \"729\u00d79=\" \"822\u00d72=\"\nReplace-Text \"825\u00d78=\" \"782\u00d77=\"\nReplace-Text \"441\u00d78=\" \"815\u00d76=\"\nReplace-Text \"858\u00d72=\" \"143\u00d72=\"\nReplace-Text \"404\u00d74=\" \"658\u00d78=\"\nReplace-Text \"737\u00d75=\" \"737\u00d78=\"\nReplace-Text \"531\u00d76=\" \"618\u00d76=\"\nReplace-Text \"660\u00d77=\" \"355\u00d72=\"\nReplace-Text \"636\u00d75=\" \"994\u00d79=\"\nReplace-Text \"910\u00d77=\" \"901\u00d73=\"\nReplace-Text \"998\u00d78=\" \"650\u00d72=\"\nReplace-Text \"198\u00d76=\" \"116\u00d78=\"\nReplace-Text \"912\u00d74=\" \"583\u00d72=\"\n"}
